$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Refresh the cached "datetimeFigureOut" footer field text
#    (slide master + every custom slide layout) from 2021/6/16 -> 2022/5/7
# ---------------------------------------------------------------------------
function Update-DateShape($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $sh = $container.Shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "2021/6/16") {
                $sh.TextFrame.TextRange.Text = "2022/5/7"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShape($master)
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    Update-DateShape($master.CustomLayouts.Item($li))
}

# ---------------------------------------------------------------------------
# 2) "( N )" -> "( N / 4 )" page-counter textboxes on the song slides.
#    Each number (1..4) appears on a pair of slides: the first slide of the
#    pair ends up with the text split across three runs, the second slide
#    simply gets one run with the full replacement text.
# ---------------------------------------------------------------------------

# slide 11: "( 4 )" -> three runs "( " / "4 " / "/ 4 )"
$sh = $p.Slides.Item(11).Shapes.Item("TextBox 2")
$tr = $sh.TextFrame.TextRange
$tr.Text = "( "
$tr.InsertAfter("4 ") | Out-Null
$tr.InsertAfter("/ 4 )") | Out-Null

# slide 12: "( 4 )" -> single run "( 4 / 4 )"
$p.Slides.Item(12).Shapes.Item("TextBox 2").TextFrame.TextRange.Text = "( 4 / 4 )"

# slide 2: "( 1 )" -> three runs "( " / "1 / 4 " / ")"
$sh = $p.Slides.Item(2).Shapes.Item("TextBox 2")
$tr = $sh.TextFrame.TextRange
$tr.Text = "( "
$tr.InsertAfter("1 / 4 ") | Out-Null
$tr.InsertAfter(")") | Out-Null

# slide 3: "( 1 )" -> single run "( 1 / 4 )"
$p.Slides.Item(3).Shapes.Item("TextBox 2").TextFrame.TextRange.Text = "( 1 / 4 )"

# slide 5: "( 2 )" -> three runs "( " / "2 " / "/ 4 )"
$sh = $p.Slides.Item(5).Shapes.Item("TextBox 2")
$tr = $sh.TextFrame.TextRange
$tr.Text = "( "
$tr.InsertAfter("2 ") | Out-Null
$tr.InsertAfter("/ 4 )") | Out-Null

# slide 6: "( 2 )" -> single run "( 2 / 4 )"
$p.Slides.Item(6).Shapes.Item("TextBox 2").TextFrame.TextRange.Text = "( 2 / 4 )"

# slide 8: "( 3 )" -> three runs "( " / "3 " / "/ 4 )"
$sh = $p.Slides.Item(8).Shapes.Item("TextBox 2")
$tr = $sh.TextFrame.TextRange
$tr.Text = "( "
$tr.InsertAfter("3 ") | Out-Null
$tr.InsertAfter("/ 4 )") | Out-Null

# slide 9: "( 3 )" -> single run "( 3 / 4 )"
$p.Slides.Item(9).Shapes.Item("TextBox 2").TextFrame.TextRange.Text = "( 3 / 4 )"
